$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update previous_count (C) and change (D) values for rows where the
# previous count was refreshed to match the current count (change -> 0),
# reflecting the new archived snapshot file used for comparison.
$ws.Range("C15").Value = 32
$ws.Range("D15").Value = 0

$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 0

$ws.Range("C27").Value = 20
$ws.Range("D27").Value = 0

$ws.Range("C29").Value = 20
$ws.Range("D29").Value = 0
